# Add a simple "keyword -> auto reply" lookup table to the first sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "keyword"
$ws.Range("B1").Value = "reply"

# Keyword / reply pairs
$ws.Range("A2").Value = "哈囉"
$ws.Range("B2").Value = "您好～請問需要什麼協助？"

$ws.Range("A3").Value = "價格"
$ws.Range("B3").Value = "價格請參考賣貨便或蝦皮連結"

$ws.Range("A4").Value = "出貨"
$ws.Range("B4").Value = "最晚明天出貨喔"

# Match the author's last selection when they saved the file
$ws.Range("E18").Select()
